# tilda6-BOM.xlsx edit:
# - Swap undervoltage monitor part (ME2807A33M3G -> SSP61CC3002MR) used for U8,
#   and un-DNP the matching footprint on the top board (SW2/SW3 row), now
#   priced and linked to the new part's LCSC page.
# - R12 position now also covers R17 (qty 1 -> 2).
# - Update a couple of hyperlink display texts so they show the full URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM costing")

$newUrl = "https://lcsc.com/product-detail/Monitors-Reset-Circuits_Shanghai-Siproin-Microelectronics-SSP61CC3002MR_C277924.html"

# --- Row 13: "> SW2, SW3" (Button_Switch_SMD:SW_SPST_B3U-1000P) ---
# was marked DNP (price 0); now populated with the SSP61CC3002MR undervoltage
# monitor footprint/price/link.
$ws.Range("F13").ClearContents()
$ws.Range("K13").Value = 0.016
$ws.Range("N13").Value = $newUrl

# --- Row 40: "    R12" resistor, now shared with R17 ---
$ws.Range("A40").Value = "    R12, R17"
$ws.Range("E40").Value = 2

# --- Row 53: "    U8" ME2807A33M3G -> SSP61CC3002MR ---
$ws.Range("B53").Value = "SSP61CC3002MR"
$ws.Range("K53").Value = 0.04
$ws.Range("M53").Value = 225000
$ws.Range("N53").Value = $newUrl

# --- Hyperlink display text fixups ---
foreach ($h in $ws.Hyperlinks) {
    if ($h.TextToDisplay -eq "https://lcsc.com/product-detail/Multi-Directional-Switches_Korean-Hroparts-Elec-K1-5202UA-02_C145900.html") {
        $h.TextToDisplay = "alt: https://lcsc.com/product-detail/Multi-Directional-Switches_Korean-Hroparts-Elec-K1-5202UA-02_C145900.html"
    }
    if ($h.TextToDisplay -eq "https://aliexpress") {
        $h.TextToDisplay = "https://aliexpress.com/item/32968351207.html"
    }
}

# --- Reflect the final selection used while making the E40 edit ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E40").Select()
